$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -1
